$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.108500653181125
$ws.Range("D2").Value = 0.1624034547986639
$ws.Range("E2").Value = 0.2237722243905651
$ws.Range("F2").Value = 1.894015224933568
$ws.Range("G2").Value = 1.37122743773611
$ws.Range("H2").Value = 1.244373562990262
$ws.Range("J2").Value = 0.3312715687186483
$ws.Range("L2").Value = 0.9779554647491295
$ws.Range("B3").Value = 1.986908029744598
$ws.Range("D3").Value = 0.155973693295536
$ws.Range("E3").Value = 0.2133343387172886
$ws.Range("F3").Value = 1.910661287167287
$ws.Range("G3").Value = 1.366833221057931
$ws.Range("H3").Value = 1.251570726095522
$ws.Range("J3").Value = 0.3143953838168159
$ws.Range("L3").Value = 0.8795472021727733
$ws.Range("B4").Value = 1.91263189187049
$ws.Range("D4").Value = 0.1519987860202292
$ws.Range("E4").Value = 0.2069256154009835
$ws.Range("F4").Value = 1.922695372693141
$ws.Range("G4").Value = 1.365679369274773
$ws.Range("H4").Value = 1.257040276602012
$ws.Range("J4").Value = 0.3040608397307238
$ws.Range("L4").Value = 0.8190558897381663
$ws.Range("B5").Value = 1.882460588766492
$ws.Range("D5").Value = 0.1503721614365219
$ws.Range("E5").Value = 0.2043142071994737
$ws.Range("F5").Value = 1.928053358282369
$ws.Range("G5").Value = 1.365594887783487
$ws.Range("H5").Value = 1.259532287943372
$ws.Range("J5").Value = 0.2998566590362088
$ws.Range("L5").Value = 0.794388677716654
$ws.Range("B6").Value = 1.877456537124829
$ws.Range("D6").Value = 0.1501016491847196
$ws.Range("E6").Value = 0.203880601471333
$ws.Range("F6").Value = 1.928970410208308
$ws.Range("G6").Value = 1.365604080936762
$ws.Range("H6").Value = 1.259961946795997
$ws.Range("J6").Value = 0.2991590041704484
$ws.Range("L6").Value = 0.7902917288330684
$ws.Range("B7").Value = 1.912224598113085
$ws.Range("D7").Value = 0.1519768764261045
$ws.Range("E7").Value = 0.2068903960234749
$ws.Range("F7").Value = 1.922765796786393
$ws.Range("G7").Value = 1.365676671583074
$ws.Range("H7").Value = 1.257072820727217
$ws.Range("J7").Value = 0.3040041109791218
$ws.Range("L7").Value = 0.8187232850378336
$ws.Range("B8").Value = 2.066496556876587
$ws.Range("D8").Value = 0.1601920601585363
$ws.Range("E8").Value = 0.220173290808269
$ws.Range("F8").Value = 1.899377234305945
$ws.Range("G8").Value = 1.369390247456892
$ws.Range("H8").Value = 1.246636502507727
$ws.Range("J8").Value = 0.3254471320186099
$ws.Range("L8").Value = 0.9440388041579695
$ws.Range("B9").Value = 2.372043893349087
$ws.Range("D9").Value = 0.1760897513465949
$ws.Range("E9").Value = 0.2462172810919583
$ws.Range("F9").Value = 1.867990958150244
$ws.Range("G9").Value = 1.389044665581878
$ws.Range("H9").Value = 1.234553223783422
$ws.Range("J9").Value = 0.3677039152485406
$ws.Range("L9").Value = 1.189226715701352
$ws.Range("B10").Value = 2.59838254218198
$ws.Range("D10").Value = 0.1876444298038535
$ws.Range("E10").Value = 0.2653449087005555
$ws.Range("F10").Value = 1.85388540771946
$ws.Range("G10").Value = 1.411196812478209
$ws.Range("H10").Value = 1.230853100196327
$ws.Range("J10").Value = 0.3988652137858679
$ws.Range("L10").Value = 1.369030366271204
$ws.Range("B11").Value = 2.701757446621116
$ws.Range("D11").Value = 0.1928747170071432
$ws.Range("E11").Value = 0.2740443249909106
$ws.Range("F11").Value = 1.849440436939659
$ws.Range("G11").Value = 1.422986010767829
$ws.Range("H11").Value = 1.230308815569572
$ws.Range("J11").Value = 0.4130644161242856
$ws.Range("L11").Value = 1.450757496733672
$ws.Range("B12").Value = 2.740962091939593
$ws.Range("D12").Value = 0.1948516119541068
$ws.Range("E12").Value = 0.27733820599898
$ws.Range("F12").Value = 1.848043029040838
$ws.Range("G12").Value = 1.427699422124704
$ws.Range("H12").Value = 1.230267675145768
$ws.Range("J12").Value = 0.4184444881133231
$ws.Range("L12").Value = 1.481695757096929
$ws.Range("B13").Value = 2.732516061359604
$ws.Range("D13").Value = 0.1944260162613034
$ws.Range("E13").Value = 0.2766288296241157
$ws.Range("F13").Value = 1.848331238119783
$ws.Range("G13").Value = 1.426673179513045
$ws.Range("H13").Value = 1.230269179345981
$ws.Range("J13").Value = 0.417285657867609
$ws.Range("L13").Value = 1.475033102279383
$ws.Range("B14").Value = 2.704981660217186
$ws.Range("D14").Value = 0.1930374313039067
$ws.Range("E14").Value = 0.2743153233140063
$ws.Range("F14").Value = 1.849319730281621
$ws.Range("G14").Value = 1.423368776964281
$ws.Range("H14").Value = 1.230302117841291
$ws.Range("J14").Value = 0.4135069764486445
$ws.Range("L14").Value = 1.453303009588524
$ws.Range("B15").Value = 2.688123690460316
$ws.Range("D15").Value = 0.1921864020392263
$ws.Range("E15").Value = 0.272898178019318
$ws.Range("F15").Value = 1.849962498465459
$ws.Range("G15").Value = 1.421377261388727
$ws.Range("H15").Value = 1.23034381295227
$ws.Range("J15").Value = 0.4111928261926039
$ws.Range("L15").Value = 1.439991380938068
$ws.Range("B16").Value = 2.591635003663384
$ws.Range("D16").Value = 0.1873020982342695
$ws.Range("E16").Value = 0.2647763327395509
$ws.Range("F16").Value = 1.854215767154329
$ws.Range("G16").Value = 1.410461059727851
$ws.Range("H16").Value = 1.230911690207819
$ws.Range("J16").Value = 0.3979377169908105
$ws.Range("L16").Value = 1.36368792717019
$ws.Range("B17").Value = 2.532547564287654
$ws.Range("D17").Value = 0.184299103527664
$ws.Range("E17").Value = 0.2597932751706296
$ws.Range("F17").Value = 1.857331627680963
$ws.Range("G17").Value = 1.404204886871412
$ws.Range("H17").Value = 1.23155261486491
$ws.Range("J17").Value = 0.389812025062696
$ws.Range("L17").Value = 1.316860855888308
$ws.Range("B18").Value = 2.498600865977721
$ws.Range("D18").Value = 0.1825694113561269
$ws.Range("E18").Value = 0.2569269869784492
$ws.Range("F18").Value = 1.859309238466821
$ws.Range("G18").Value = 1.400767475666072
$ws.Range("H18").Value = 1.232028386422314
$ws.Range("J18").Value = 0.3851405952023725
$ws.Range("L18").Value = 1.289920933335907
$ws.Range("B19").Value = 2.487113782569395
$ws.Range("D19").Value = 0.1819833456139719
$ws.Range("E19").Value = 0.255956486873103
$ws.Range("F19").Value = 1.860010605630265
$ws.Range("G19").Value = 1.399631186559191
$ws.Range("H19").Value = 1.232207837744852
$ws.Range("J19").Value = 0.3835593249174565
$ws.Range("L19").Value = 1.280798483808553
$ws.Range("B20").Value = 2.538833507565869
$ws.Range("D20").Value = 0.1846190306501967
$ws.Range("E20").Value = 0.2603237482574841
$ws.Range("F20").Value = 1.856980729184258
$ws.Range("G20").Value = 1.404854185995219
$ws.Range("H20").Value = 1.231473291327148
$ws.Range("J20").Value = 0.3906767876770658
$ws.Range("L20").Value = 1.321846325913043
$ws.Range("B21").Value = 2.713067588626473
$ws.Range("D21").Value = 0.1934453925121602
$ws.Range("E21").Value = 0.2749948680749128
$ws.Range("F21").Value = 1.849021611261193
$ws.Range("G21").Value = 1.424332576410478
$ws.Range("H21").Value = 1.230287956252027
$ws.Range("J21").Value = 0.414616783640696
$ws.Range("L21").Value = 1.459685940301824
$ws.Range("B22").Value = 2.82728225887297
$ws.Range("D22").Value = 0.1991923751856604
$ws.Range("E22").Value = 0.2845809082402653
$ws.Range("F22").Value = 1.845486831532142
$ws.Range("G22").Value = 1.438516125160987
$ws.Range("H22").Value = 1.230475394260338
$ws.Range("J22").Value = 0.4302811504855697
$ws.Range("L22").Value = 1.549713448976888
$ws.Range("B23").Value = 2.766292531599618
$ws.Range("D23").Value = 0.1961270618392916
$ws.Range("E23").Value = 0.27946491938701
$ws.Range("F23").Value = 1.847220126579117
$ws.Range("G23").Value = 1.430812157106686
$ws.Range("H23").Value = 1.230286923790004
$ws.Range("J23").Value = 0.4219192063569892
$ws.Range("L23").Value = 1.50166958299252
$ws.Range("B24").Value = 2.53599155950667
$ws.Range("D24").Value = 0.1844744016637918
$ws.Range("E24").Value = 0.26008392591838
$ws.Range("F24").Value = 1.857138790405116
$ws.Range("G24").Value = 1.404560141753961
$ws.Range("H24").Value = 1.231508819387585
$ws.Range("J24").Value = 0.3902858280462311
$ws.Range("L24").Value = 1.319592452237941
$ws.Range("B25").Value = 2.289060904840881
$ws.Range("D25").Value = 0.1718113205381684
$ws.Range("E25").Value = 0.239172591442653
$ws.Range("F25").Value = 1.874918619304935
$ws.Range("G25").Value = 1.382385626093885
$ws.Range("H25").Value = 1.236917754099238
$ws.Range("J25").Value = 0.3562513486015746
$ws.Range("L25").Value = 1.122956266358358
